$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item(1)

# Match the final selection state expected on the original sheet
# (a "select all cells" selection, sqref spanning the whole grid).
$src.Range("A1:XFD1048576").Select()

# Duplicate "Attached Functionality" -> new sheet lands right after it.
$src.Copy($null, $src)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Attached Functionality xBB"

# Make room for a new "Black Box" row just above the trailing
# "Wg" / "Attached Functionality" rows, preserving their formatting.
$vLast = $ws2.Range("A14").Value()
$vPrev = $ws2.Range("A13").Value()

$ws2.Range("A14").Copy()
$ws2.Range("A15").PasteSpecial(-4122)

$ws2.Range("A15").Value = $vLast
$ws2.Range("A14").Value = $vPrev
$ws2.Range("A13").Value = "Black Box"

# New sheet becomes the active tab, selection resting on A11.
$ws2.Range("A11").Select()
$ws2.Activate()
